# Add a new "2022" column (N) to the table, mirroring the existing 2021
# column (M): same formatting per row, one new data point per metric row,
# then move the active selection to N15 (just below the new data), matching
# the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: header year value
$ws.Range("M4").Copy($ws.Range("N4")) | Out-Null
$ws.Range("N4").Value = 2022

# Row 5
$ws.Range("M5").Copy($ws.Range("N5")) | Out-Null
$ws.Range("N5").Value = 4.3

# Row 6
$ws.Range("M6").Copy($ws.Range("N6")) | Out-Null
$ws.Range("N6").Value = 5.0999999999999996

# Row 7
$ws.Range("M7").Copy($ws.Range("N7")) | Out-Null
$ws.Range("N7").Value = 3.1

# Row 8
$ws.Range("M8").Copy($ws.Range("N8")) | Out-Null
$ws.Range("N8").Value = 2.9

# Row 9
$ws.Range("M9").Copy($ws.Range("N9")) | Out-Null
$ws.Range("N9").Value = 3.4

# Row 10
$ws.Range("M10").Copy($ws.Range("N10")) | Out-Null
$ws.Range("N10").Value = 2.2999999999999998

# Row 11
$ws.Range("M11").Copy($ws.Range("N11")) | Out-Null
$ws.Range("N11").Value = 92.8

# Row 12
$ws.Range("M12").Copy($ws.Range("N12")) | Out-Null
$ws.Range("N12").Value = 91.6

# Row 13
$ws.Range("M13").Copy($ws.Range("N13")) | Out-Null
$ws.Range("N13").Value = 94.6

# Match the final selection recorded in the saved workbook.
$ws.Range("N15").Select() | Out-Null
